# Hortaliza, Vega Central Mapocho de Santiago - Cebollín
# Weekly data refresh: a new reporting date (44505) is inserted at the
# top of the data block (rows 575-576), pushing all the existing rows
# down by two. This grows the sheet from A1:R670 to A1:R672.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the current first data row of this
# block (row 575); everything from 575 downward shifts to 577 downward.
$ws.Rows("575:576").Insert()

# --- New row 575: Cebollín "Primera", fecha 44505 ---
$ws.Cells.Item(575, 1).Value2  = 9
$ws.Cells.Item(575, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(575, 3).Value2  = "Metropolitana"
$ws.Cells.Item(575, 4).Value2  = 44505
$ws.Cells.Item(575, 5).Value2  = 13
$ws.Cells.Item(575, 6).Value2  = 100112037
$ws.Cells.Item(575, 7).Value2  = "Cebollín"
$ws.Cells.Item(575, 8).Value2  = "Sin especificar"
$ws.Cells.Item(575, 9).Value2  = "Primera"
$ws.Cells.Item(575, 10).Value2 = 250
$ws.Cells.Item(575, 11).Value2 = 2300
$ws.Cells.Item(575, 12).Value2 = 2500
$ws.Cells.Item(575, 13).Value2 = 2400
$ws.Cells.Item(575, 14).Value2 = "$/paquete 36 unidades"
$ws.Cells.Item(575, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(575, 16).Value2 = 67
$ws.Cells.Item(575, 17).Value2 = 36
$ws.Cells.Item(575, 18).Value2 = "Hortaliza"

# --- New row 576: Cebollín "Segunda", fecha 44505 ---
$ws.Cells.Item(576, 1).Value2  = 9
$ws.Cells.Item(576, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(576, 3).Value2  = "Metropolitana"
$ws.Cells.Item(576, 4).Value2  = 44505
$ws.Cells.Item(576, 5).Value2  = 13
$ws.Cells.Item(576, 6).Value2  = 100112037
$ws.Cells.Item(576, 7).Value2  = "Cebollín"
$ws.Cells.Item(576, 8).Value2  = "Sin especificar"
$ws.Cells.Item(576, 9).Value2  = "Segunda"
$ws.Cells.Item(576, 10).Value2 = 160
$ws.Cells.Item(576, 11).Value2 = 1800
$ws.Cells.Item(576, 12).Value2 = 2000
$ws.Cells.Item(576, 13).Value2 = 1900
$ws.Cells.Item(576, 14).Value2 = "$/paquete 36 unidades"
$ws.Cells.Item(576, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(576, 16).Value2 = 53
$ws.Cells.Item(576, 17).Value2 = 36
$ws.Cells.Item(576, 18).Value2 = "Hortaliza"
